# Generate Report for Archive
# The localization status of file "2b292678-d178-4b57-924f-dc2bc3f548a8.md"
# moved from "Ready for handoff" to "In Translation" for both the zh-cn and
# de-de locales. Update the Overview summary sheet plus each per-locale
# detail sheet to reflect the new status.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "In Translation"
$overview.Range("C2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "In Translation"
